# Insert a new weekly price record for "Ciboulette" at row 276, shifting the
# existing rows 276:291 down to 277:292 (dimension grows from R291 to R292).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down by one row, starting at row 276.
$ws.Rows("276:276").Insert()

$row = 276

$ws.Cells.Item($row, 1).Value = 9
$ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($row, 3).Value = "Metropolitana"
$ws.Cells.Item($row, 4).Value = 44516
$ws.Cells.Item($row, 5).Value = 13
$ws.Cells.Item($row, 6).Value = 100112039
$ws.Cells.Item($row, 7).Value = "Ciboulette"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 250
$ws.Cells.Item($row, 11).Value = 1000
$ws.Cells.Item($row, 12).Value = 1200
$ws.Cells.Item($row, 13).Value = 1100
$ws.Cells.Item($row, 14).Value = "`$/docena de atados"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 367
$ws.Cells.Item($row, 17).Value = 3
$ws.Cells.Item($row, 18).Value = "Hortaliza"
